$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: shift labels left (max -> E1, prediction -> C1, rejection-f -> D1)
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Data rows: C becomes the taxonomy string (was in D), D stays the taxonomy
# string, E becomes the numeric flag 1 (was the taxonomy string)
$ws.Range("C2").Value = "g__QAMX01"
$ws.Range("D2").Value = "g__QAMX01"
$ws.Range("E2").Value = 1

$ws.Range("C3").Value = "g__QAMX01"
$ws.Range("D3").Value = "g__QAMX01"
$ws.Range("E3").Value = 1

$ws.Range("C4").Value = "g__QAMX01"
$ws.Range("D4").Value = "g__QAMX01"
$ws.Range("E4").Value = 1

$ws.Range("C5").Value = "g__QAMX01"
$ws.Range("D5").Value = "g__QAMX01"
$ws.Range("E5").Value = 1

$ws.Range("C6").Value = "g__QAMX01"
$ws.Range("D6").Value = "g__QAMX01"
$ws.Range("E6").Value = 1
